$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2-30: Community, Allocated Population, Shelter Assigned, Level
$data = @(
    @(2, 'Balite', 602, 'San Marcos National H.S.', 1),
    @(3, 'Balungao', 687, 'San Marcos National H.S.', 1),
    @(4, 'Buguion', 461, 'San Marcos National H.S.', 1),
    @(5, 'Bulusan', 313, 'F. Mendoza Memorial Elem Sch.', 1),
    @(6, 'Calizon', 267, 'BMLTC Multi-Purpose Bldg and EC', 1),
    @(7, 'Calumpang', 423, 'San Marcos National H.S.', 1),
    @(8, 'Caniogan', 542, 'F. Mendoza Memorial Elem Sch.', 1),
    @(9, 'Corazon', 261, 'BMLTC Multi-Purpose Bldg and EC', 1),
    @(10, 'Frances', 736, 'Palimbang Empty Lot', 1),
    @(11, 'Gatbuca', 767, 'NV9 Multi-Purpose', 1),
    @(12, 'Gugo', 236, 'Palimbang Empty Lot', 1),
    @(13, 'Iba Este', 500, 'NV9 Multi-Purpose', 1),
    @(14, 'Iba O''Este', 1691, 'San Marcos National H.S.', 1),
    @(15, 'Longos', 512, 'San Marcos Elem. Sch.', 2),
    @(16, 'Meysulao', 514, 'BMLTC Multi-Purpose Bldg and EC', 1),
    @(17, 'Meyto', 351, 'Calumpit Sports Complex', 1),
    @(18, 'Palimbang', 203, 'Palimbang Empty Lot', 1),
    @(19, 'Panducot', 211, 'BMLTC Multi-Purpose Bldg and EC', 1),
    @(20, 'Pio Cruzcosa', 560, 'Palimbang Empty Lot', 1),
    @(21, 'Poblacion', 215, 'Calumpit Sports Complex', 1),
    @(22, 'Pungo', 1144, 'Doña Damiana Elem School', 1),
    @(23, 'San Jose', 680, 'San Marcos Elem. Sch.', 2),
    @(24, 'San Marcos', 321, 'San Marcos Elem. Sch.', 2),
    @(25, 'San Miguel', 721, 'Doña Damiana Elem School', 1),
    @(26, 'Santa Lucia', 296, 'BMLTC Multi-Purpose Bldg and EC', 1),
    @(27, 'Santo Niño', 306, 'F. Mendoza Memorial Elem Sch.', 1),
    @(28, 'Sapang Bayan', 377, 'San Marcos Elem. Sch.', 2),
    @(29, 'Sergio Bayan', 208, 'F. Mendoza Memorial Elem Sch.', 1),
    @(30, 'Sucol', 128, 'F. Mendoza Memorial Elem Sch.', 1)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
}

Write-Output "Done updating rows"
